$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SiteDevelopmentValues")

# Insert a new row above row 18, shifting the existing row 18 (and below) down
$ws.Rows.Item(18).Insert()

# Fill the new row 18 with the site-rental entry
$ws.Range("A18").Value = "existing_site_rental_per_night"
$ws.Range("B18").Value = "Cost to use a telescope per night if we don't own it"
$ws.Range("C18").Value = 10000
$ws.Range("D18").Value = "This is what it costs for a block of time from SMA?"

# Match the currency/accounting number format used elsewhere in column C
$ws.Range("C18").NumberFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

# Match the wrapped-text row height used by the rest of this sheet
$ws.Rows.Item(18).RowHeight = 36

$ws.Range("C18").Select() | Out-Null
